# Add season record columns (Wins, Losses, Ties) to the HOU 1999 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns AD/AE/AF, formatted like the rest of
# the header row (bold, thin border, centered / top-aligned) ---
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1         # xlContinuous (thin, all sides)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-40: season record is the same for every player row ---
$lastRow = 40
$ws.Range("AD2:AD$lastRow").Value = 97
$ws.Range("AE2:AE$lastRow").Value = 65
$ws.Range("AF2:AF$lastRow").Value = 0

Write-Output "done"
